# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# All cells in columns B:E are stored as text in this sheet, so any
# replacement value that Excel would otherwise auto-convert to a number
# (plain decimals like "206.14") is written with a leading apostrophe to
# force text, matching the original cell typing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.739.55"
$ws.Range("E2").Value = "  -2.40%  "
$ws.Range("D3").Value = "1.560.51"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'206.14"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'21.91"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").Value = "'0.0861"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "1.782.02"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "1.563.01"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").Value = "'61.49"
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("D17").Value = "26.744.38"
$ws.Range("E17").Value = "  -2.38%  "
$ws.Range("D18").Value = "'7.35"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'4.09"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "'9.33"
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").Value = "'152.51"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").Value = "'6.76"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").Value = "'14.82"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("E30").Value = "  -4.13%  "
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("D33").Value = "1.382.70"
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("D34").Value = "'2.92"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("E37").Value = "  -3.62%  "
$ws.Range("D38").Value = "'0.0163"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").Value = "'0.520"
$ws.Range("E39").Value = "  -2.26%  "
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D42").Value = "'0.993"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").Value = "'1.76"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("D46").Value = "'63.13"
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("D47").Value = "1.695.53"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").Value = "'85.37"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").Value = "0.0₇0982"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0948"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0493"
$ws.Range("E51").Value = "  -0.32%  "
